$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ28302677",
    "summ28549760",
    "summ28819365",
    "summ29086549",
    "summ29348594",
    "summ29614986",
    "summ29879675",
    "summ30140254",
    "summ30398302"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
